$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (RS order = 1)
$ws.Range("E2").Value = 0.8825501377182884
$ws.Range("I2").Value = 235.2063318150511
$ws.Range("J2").Value = 10.22636225282831

# Row 3 (RS order = 2)
$ws.Range("E3").Value = 0.1251035512923761
$ws.Range("I3").Value = 382.493470473173
$ws.Range("J3").Value = 127.4978234910577

# Row 4 (RS order = 3)
$ws.Range("E4").Value = 0.05987834989280399
$ws.Range("I4").Value = 220.7374005464412
$ws.Range("J4").Value = 220.7374005464412
